# Apply crypto price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.568.76'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '3.040.02'
$ws.Range('E3').Value = '  -3.01%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.21'
$ws.Range('E5').Value = '  -2.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '612.63'
$ws.Range('E6').Value = '  -3.78%  '
$ws.Range('E7').Value = '  -9.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.874'
$ws.Range('E8').Value = '  +14.71%  '
$ws.Range('D10').Value = '3.036.92'
$ws.Range('E10').Value = '  -3.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.669'
$ws.Range('E11').Value = '  +19.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.188'
$ws.Range('E12').Value = '  +4.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000239'
$ws.Range('E13').Value = '  -5.15%  '
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').Value = '89.293.16'
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.29'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '3.607.52'
$ws.Range('E17').Value = '  -2.68%  '
$ws.Range('D18').Value = '3.043.96'
$ws.Range('E18').Value = '  -2.98%  '
$ws.Range('E19').Value = '  -1.99%  '
$ws.Range('E20').Value = '  -5.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.37'
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '423.56'
$ws.Range('E22').Value = '  -1.07%  '
$ws.Range('E23').Value = '  -1.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.01'
$ws.Range('E24').Value = '  +1.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.37'
$ws.Range('E25').Value = '  -1.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '83.25'
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.59'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.161'
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.30'
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.71'
$ws.Range('E32').Value = '  -7.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '499.55'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.64'
$ws.Range('E34').Value = '  -6.85%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.80'
$ws.Range('E35').Value = '  -1.98%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.71'
$ws.Range('E36').Value = '  +3.00%  '
$ws.Range('E37').Value = '  -4.44%  '
$ws.Range('E38').Value = '  -10.18%  '
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.361'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.83'
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.135'
$ws.Range('E44').Value = '  +4.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.13'
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0686'
$ws.Range('E47').Value = '  +7.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.16'
$ws.Range('E48').Value = '  +5.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '162.63'
$ws.Range('E49').Value = '  -1.29%  '
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.589'
$ws.Range('E51').Value = '  -1.51%  '
